# Apply the corrected/shifted naive-forecaster YoY data (bugfix: drop stale row, recompute C/E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = 39583
$ws.Range("B2").Value2 = 2008
$ws.Range("D2").Value2 = 2009
$ws.Range("E2").Value2 = 1.328558632615739
# Row 3
$ws.Range("A3").Value2 = 39765
$ws.Range("B3").Value2 = 2008
$ws.Range("D3").Value2 = 2009
$ws.Range("E3").Value2 = 0.4282194198276246
# Row 4
$ws.Range("A4").Value2 = 39948
$ws.Range("B4").Value2 = 2009
$ws.Range("C4").Value2 = -1.435981453719049
$ws.Range("D4").Value2 = 2010
$ws.Range("E4").Value2 = -0.7704417043119083
# Row 5
$ws.Range("A5").Value2 = 40130
$ws.Range("B5").Value2 = 2009
$ws.Range("C5").Value2 = -1.61188520630966
$ws.Range("D5").Value2 = 2010
$ws.Range("E5").Value2 = -0.8235211753995442
# Row 6
$ws.Range("A6").Value2 = 40310
$ws.Range("B6").Value2 = 2010
$ws.Range("C6").Value2 = 0.406633294022174
$ws.Range("D6").Value2 = 2011
$ws.Range("E6").Value2 = -0.1034614224434405
# Row 7
$ws.Range("A7").Value2 = 40494
$ws.Range("B7").Value2 = 2010
$ws.Range("C7").Value2 = 1.128600547465064
$ws.Range("D7").Value2 = 2011
$ws.Range("E7").Value2 = 1.121293995080253
# Row 8
$ws.Range("A8").Value2 = 40676
$ws.Range("B8").Value2 = 2011
$ws.Range("C8").Value2 = 2.185496833134781
$ws.Range("D8").Value2 = 2012
$ws.Range("E8").Value2 = 0.6652762968575532
# Row 9
$ws.Range("A9").Value2 = 40862
$ws.Range("B9").Value2 = 2011
$ws.Range("C9").Value2 = 2.397632938760519
$ws.Range("D9").Value2 = 2012
$ws.Range("E9").Value2 = 1.665250327443002
# Row 10
$ws.Range("A10").Value2 = 41044
$ws.Range("B10").Value2 = 2012
$ws.Range("C10").Value2 = 0.8574941660507873
$ws.Range("D10").Value2 = 2013
$ws.Range("E10").Value2 = 1.693469135756587
# Row 11
$ws.Range("A11").Value2 = 41228
$ws.Range("B11").Value2 = 2012
$ws.Range("C11").Value2 = 0.5991205513815823
$ws.Range("D11").Value2 = 2013
$ws.Range("E11").Value2 = 1.079796209653616
# Row 12
$ws.Range("A12").Value2 = 41409
$ws.Range("B12").Value2 = 2013
$ws.Range("C12").Value2 = 0.1494732105682406
$ws.Range("D12").Value2 = 2014
$ws.Range("E12").Value2 = 0.8024032015999882
# Row 13
$ws.Range("A13").Value2 = 41592
$ws.Range("B13").Value2 = 2013
$ws.Range("C13").Value2 = 0.5555179840670776
$ws.Range("D13").Value2 = 2014
$ws.Range("E13").Value2 = 1.374377011838535
# Row 14
$ws.Range("A14").Value2 = 41774
$ws.Range("B14").Value2 = 2014
$ws.Range("C14").Value2 = 1.656936590801972
$ws.Range("D14").Value2 = 2015
$ws.Range("E14").Value2 = 0.922773818606859
# Row 15
$ws.Range("A15").Value2 = 41957
$ws.Range("B15").Value2 = 2014
$ws.Range("C15").Value2 = 1.749602965204744
$ws.Range("D15").Value2 = 2015
$ws.Range("E15").Value2 = 1.310895847186577
# Row 16
$ws.Range("A16").Value2 = 42137
$ws.Range("B16").Value2 = 2015
$ws.Range("C16").Value2 = 1.346932828201242
$ws.Range("D16").Value2 = 2016
$ws.Range("E16").Value2 = 1.364302026343633
# Row 17
$ws.Range("A17").Value2 = 42321
$ws.Range("B17").Value2 = 2015
$ws.Range("C17").Value2 = 1.513781691628258
$ws.Range("D17").Value2 = 2016
$ws.Range("E17").Value2 = 1.862478303083726
# Row 18
$ws.Range("A18").Value2 = 42503
$ws.Range("B18").Value2 = 2016
$ws.Range("C18").Value2 = 1.745747589686109
$ws.Range("D18").Value2 = 2017
$ws.Range("E18").Value2 = 1.644798626926303
# Row 19
$ws.Range("A19").Value2 = 42689
$ws.Range("B19").Value2 = 2016
$ws.Range("C19").Value2 = 1.72540577912379
$ws.Range("D19").Value2 = 2017
$ws.Range("E19").Value2 = 1.639776099317536
# Row 20
$ws.Range("A20").Value2 = 42867
$ws.Range("B20").Value2 = 2017
$ws.Range("C20").Value2 = 1.843649045891893
$ws.Range("D20").Value2 = 2018
$ws.Range("E20").Value2 = 1.741128155516525
# Row 21
$ws.Range("A21").Value2 = 43053
$ws.Range("B21").Value2 = 2017
$ws.Range("C21").Value2 = 2.026192376700298
$ws.Range("D21").Value2 = 2018
$ws.Range("E21").Value2 = 2.181728312936415
# Row 22
$ws.Range("A22").Value2 = 43145
$ws.Range("B22").Value2 = 2018
$ws.Range("C22").Value2 = 2.284406789710336
$ws.Range("D22").Value2 = 2019
$ws.Range("E22").Value2 = 1.990690441067144
# Row 23
$ws.Range("A23").Value2 = 43235
$ws.Range("B23").Value2 = 2018
$ws.Range("C23").Value2 = 2.463589365374652
$ws.Range("D23").Value2 = 2019
$ws.Range("E23").Value2 = 2.149194501693219
# Row 24
$ws.Range("A24").Value2 = 43326
$ws.Range("B24").Value2 = 2018
$ws.Range("C24").Value2 = 2.349806433215029
$ws.Range("D24").Value2 = 2019
$ws.Range("E24").Value2 = 2.036910005299108
# Row 25
$ws.Range("A25").Value2 = 43418
$ws.Range("B25").Value2 = 2018
$ws.Range("C25").Value2 = 2.344166347125687
$ws.Range("D25").Value2 = 2019
$ws.Range("E25").Value2 = 2.010025322622599
# Row 26
$ws.Range("A26").Value2 = 43510
$ws.Range("B26").Value2 = 2019
$ws.Range("C26").Value2 = 1.665971362160357
$ws.Range("D26").Value2 = 2020
$ws.Range("E26").Value2 = 2.031292234149706
# Row 27
$ws.Range("A27").Value2 = 43600
$ws.Range("B27").Value2 = 2019
$ws.Range("C27").Value2 = 1.332860091726285
$ws.Range("D27").Value2 = 2020
$ws.Range("E27").Value2 = 1.799885362733189
# Row 28
$ws.Range("A28").Value2 = 43691
$ws.Range("B28").Value2 = 2019
$ws.Range("C28").Value2 = 1.029194292875912
$ws.Range("D28").Value2 = 2020
$ws.Range("E28").Value2 = 1.31420459445093
# Row 29
$ws.Range("A29").Value2 = 43783
$ws.Range("B29").Value2 = 2019
$ws.Range("C29").Value2 = 0.9005461608770915
$ws.Range("D29").Value2 = 2020
$ws.Range("E29").Value2 = 0.7771393814490102
# Row 30
$ws.Range("A30").Value2 = 43875
$ws.Range("B30").Value2 = 2020
$ws.Range("C30").Value2 = 0.4126128934655471
$ws.Range("D30").Value2 = 2021
$ws.Range("E30").Value2 = 1.156986202028509
# Row 31
$ws.Range("A31").Value2 = 43966
$ws.Range("B31").Value2 = 2020
$ws.Range("C31").Value2 = 0.2336391425753925
$ws.Range("D31").Value2 = 2021
$ws.Range("E31").Value2 = 0.9207450904090253
# Row 32
$ws.Range("A32").Value2 = 44068
$ws.Range("B32").Value2 = 2020
$ws.Range("C32").Value2 = -4.43626840667447
$ws.Range("D32").Value2 = 2021
$ws.Range("E32").Value2 = -2.63419394755392
# Row 33
$ws.Range("A33").Value2 = 44159
$ws.Range("B33").Value2 = 2020
$ws.Range("C33").Value2 = -4.43626840667447
$ws.Range("D33").Value2 = 2021
$ws.Range("E33").Value2 = -2.71887004062904
# Row 34
$ws.Range("A34").Value2 = 44251
$ws.Range("B34").Value2 = 2021
$ws.Range("C34").Value2 = -2.96879819115512
$ws.Range("D34").Value2 = 2022
$ws.Range("E34").Value2 = -2.438555173006141
# Row 35
$ws.Range("A35").Value2 = 44341
$ws.Range("B35").Value2 = 2021
$ws.Range("C35").Value2 = -2.010709456685855
$ws.Range("D35").Value2 = 2022
$ws.Range("E35").Value2 = -1.14257141002756
# Row 36
$ws.Range("A36").Value2 = 44432
$ws.Range("B36").Value2 = 2021
$ws.Range("C36").Value2 = -1.513408827666285
$ws.Range("D36").Value2 = 2022
$ws.Range("E36").Value2 = 0.7106578563214505
# Row 37
$ws.Range("A37").Value2 = 44525
$ws.Range("B37").Value2 = 2021
$ws.Range("C37").Value2 = -1.513408827666285
$ws.Range("D37").Value2 = 2022
$ws.Range("E37").Value2 = 0.4582698374457683
# Row 38
$ws.Range("A38").Value2 = 44617
$ws.Range("B38").Value2 = 2022
$ws.Range("C38").Value2 = 1.154413086110817
$ws.Range("D38").Value2 = 2023
$ws.Range("E38").Value2 = -1.166698219025086
# Row 39
$ws.Range("A39").Value2 = 44706
$ws.Range("B39").Value2 = 2022
$ws.Range("C39").Value2 = 1.5286818008164
$ws.Range("D39").Value2 = 2023
$ws.Range("E39").Value2 = -0.7118141543333012
# Row 40
$ws.Range("A40").Value2 = 44798
$ws.Range("B40").Value2 = 2022
$ws.Range("C40").Value2 = 1.618732201786743
$ws.Range("D40").Value2 = 2023
$ws.Range("E40").Value2 = -0.4630595634534385
# Row 41
$ws.Range("A41").Value2 = 44890
$ws.Range("B41").Value2 = 2022
$ws.Range("C41").Value2 = 1.618732201786743
$ws.Range("D41").Value2 = 2023
$ws.Range("E41").Value2 = 1.314675624401973
# Row 42
$ws.Range("A42").Value2 = 44981
$ws.Range("B42").Value2 = 2023
$ws.Range("C42").Value2 = 0.006126408955742235
$ws.Range("D42").Value2 = 2024
$ws.Range("E42").Value2 = 0.5370151562237302
# Row 43
$ws.Range("A43").Value2 = 45071
$ws.Range("B43").Value2 = 2023
$ws.Range("C43").Value2 = 0.001079933351455509
$ws.Range("D43").Value2 = 2024
$ws.Range("E43").Value2 = 0.6889047703476203
# Row 44
$ws.Range("A44").Value2 = 45163
$ws.Range("B44").Value2 = 2023
$ws.Range("C44").Value2 = -0.09609276733164585
$ws.Range("D44").Value2 = 2024
$ws.Range("E44").Value2 = 0.5285660612534882
# Row 45
$ws.Range("A45").Value2 = 45254
$ws.Range("B45").Value2 = 2023
$ws.Range("C45").Value2 = -0.09609276733164585
$ws.Range("D45").Value2 = 2024
$ws.Range("E45").Value2 = 0.1199358335146838
# Row 46
$ws.Range("A46").Value2 = 45345
$ws.Range("B46").Value2 = 2024
$ws.Range("C46").Value2 = -0.1634698065940632
$ws.Range("D46").Value2 = 2025
$ws.Range("E46").Value2 = -0.1145111565623136
# Row 47
$ws.Range("A47").Value2 = 45436
$ws.Range("B47").Value2 = 2024
$ws.Range("C47").Value2 = -0.00209793826797533
$ws.Range("D47").Value2 = 2025
$ws.Range("E47").Value2 = 0.286657616500996
# Row 48
$ws.Range("A48").Value2 = 45534
$ws.Range("B48").Value2 = 2024
$ws.Range("C48").Value2 = -0.02761034355766023
$ws.Range("D48").Value2 = 2025
$ws.Range("E48").Value2 = 0.08174908622293753
# Row 49
$ws.Range("A49").Value2 = 45618
$ws.Range("B49").Value2 = 2024
$ws.Range("C49").Value2 = -0.02761034355766023
$ws.Range("D49").Value2 = 2025
$ws.Range("E49").Value2 = 0.2676745853112728
# Row 50
$ws.Range("A50").Value2 = 45713
$ws.Range("B50").Value2 = 2025
$ws.Range("C50").Value2 = 0.4725905789402463
$ws.Range("D50").Value2 = 2026
$ws.Range("E50").Value2 = 0.08564335928031852
# Row 51
$ws.Range("A51").Value2 = 45800
$ws.Range("B51").Value2 = 2025
$ws.Range("C51").Value2 = 0.501314651583451
$ws.Range("D51").Value2 = 2026
$ws.Range("E51").Value2 = 0.2155158706220295
# Row 52
$ws.Range("A52").Value2 = 45891
$ws.Range("B52").Value2 = 2025
$ws.Range("C52").Value2 = 0.5152269879013183
$ws.Range("D52").Value2 = 2026
$ws.Range("E52").Value2 = 0.3332251551730891

# The source data dropped the oldest forecast row, so the sheet now has one fewer row.
$ws.Rows.Item(53).Delete()

Write-Output "Updated rows 2-52 and removed trailing row 53."
